# "edits to supp tables"
# Adds a new footnote row explaining how percent gaps were calculated,
# placed right after the existing "masking does not influence ... gaps"
# footnote, and nudges column C's width slightly narrower.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the last footnote row (row 14) so the existing
# footnote shifts down to row 15, then populate the new row with the
# new footnote text (plain/default style, matching the source sheet).
$ws.Rows.Item(14).Insert()
$ws.Cells.Item(14, 1).Value = " Percent gaps were calculated from unfiltered alignments as the total number of gaps divided by the total number of MSA positions, and represent the percentage of columns with at least one gap, averaged across all MSA replicates."
$ws.Cells.Item(14, 1).Style = "Normal"

# Slightly narrow column C.
$ws.Columns.Item(3).ColumnWidth = 13.83

# Move the active selection to the newly added footnote cell.
$null = $ws.Range("A14").Select()
